# Fix numbering in figure
#
# 1) Slide 7 ("G5" safety-case figure): the G5.x labels had a gap - the
#    boxes read G5.1, G5.2, G5.3, G5.5, G5.6, G5.7. Renumber the last
#    three down by one so the sequence is contiguous: G5.5->G5.4,
#    G5.6->G5.5, G5.7->G5.6. Only the short numbering label (the first
#    4 characters of the shape's text) is touched; the longer
#    description text that follows in the same shape is left alone.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -eq -1) {
        $tr = $sh.TextFrame.TextRange
        $label = $tr.Characters(1, 4)
        $labelText = $label.Text
        if ($labelText -eq "G5.5") {
            $label.Text = "G5.4"
        } elseif ($labelText -eq "G5.6") {
            $label.Text = "G5.5"
        } elseif ($labelText -eq "G5.7") {
            $label.Text = "G5.6"
        }
    }
}

# 2) The "updated automatically" date field cached on the slide master
#    and every slide layout was re-stamped by PowerPoint the next time
#    the deck was saved (2021-09-27 -> 2021-10-04). Refresh every
#    cached occurrence so it matches.
$oldDate = "2021-09-27"
$newDate = "2021-10-04"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -eq -1) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.HasTextFrame -eq -1) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
